$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Header cell for the new column, matching the style of the other header cells.
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$timestamps = @(
    "2021-10-05 13:42:06.619262",
    "2021-10-05 13:42:06.619274",
    "2021-10-05 13:42:06.619278",
    "2021-10-05 13:42:06.619281",
    "2021-10-05 13:42:06.619284",
    "2021-10-05 13:42:06.619288",
    "2021-10-05 13:42:06.619291",
    "2021-10-05 13:42:06.619294",
    "2021-10-05 13:42:06.619297",
    "2021-10-05 13:42:06.619300",
    "2021-10-05 13:42:06.619303",
    "2021-10-05 13:42:06.619306",
    "2021-10-05 13:42:06.619309",
    "2021-10-05 13:42:06.619312",
    "2021-10-05 13:42:06.619315",
    "2021-10-05 13:42:06.619318",
    "2021-10-05 13:42:06.619322",
    "2021-10-05 13:42:06.619325",
    "2021-10-05 13:42:06.619328",
    "2021-10-05 13:42:06.619331",
    "2021-10-05 13:42:06.619334",
    "2021-10-05 13:42:06.619337",
    "2021-10-05 13:42:06.619340",
    "2021-10-05 13:42:06.619343",
    "2021-10-05 13:42:06.619346",
    "2021-10-05 13:42:06.619349",
    "2021-10-05 13:42:06.619352",
    "2021-10-05 13:42:06.619355",
    "2021-10-05 13:42:06.619358",
    "2021-10-05 13:42:06.619361",
    "2021-10-05 13:42:06.619365",
    "2021-10-05 13:42:06.619368",
    "2021-10-05 13:42:06.619371",
    "2021-10-05 13:42:06.619374",
    "2021-10-05 13:42:06.619377",
    "2021-10-05 13:42:06.619380",
    "2021-10-05 13:42:06.619383",
    "2021-10-05 13:42:06.619386",
    "2021-10-05 13:42:06.619389",
    "2021-10-05 13:42:06.619392",
    "2021-10-05 13:42:06.619395",
    "2021-10-05 13:42:06.619399"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
